# Revert "Add some changes on fe": set the Status column (F) back to "Done"
# for every task row (18-48) in the "Sprint Info" sheet. This undoes the
# earlier edit that had changed these statuses to "Process".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Info")

for ($r = 18; $r -le 48; $r++) {
    $ws.Range("F$r").Value = "Done"
}

# Restore the cursor/selection to match the reverted state (cosmetic).
$ws.Select()
$ws.Range("J52").Select()
